# Bump the trailing "run number" suffix baked into the test-data names on
# every sheet: the "21" group becomes "22" and the "26" group becomes "27".
#
# "Sheet1" (3rd tab) derives names/emails with CONCATENATE(...,$I$2) /
# CONCATENATE(...,$I$23), so changing I2 and I23 cascades through every
# formula cell automatically. The "login" and "order" sheets instead carry
# literal (non-formula) copies of the same generated strings, so those get
# rewritten directly, cell by cell.

$wb = $excel.ActiveWorkbook

# ---- "Sheet1": bump the two numeric seeds; formulas recalc on their own ----
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Range("I2").Value = 22
$wsData.Range("I23").Value = 27

# ---- "login": literal EthanBaker21 / ...21@gmail.com -> ...22 ----
$namesGroup1 = @(
    "EthanBaker", "DelanieCarman", "BretAgnew", "EdgardoTaylor", "TyrekReis",
    "LeannaChow", "TuckerCarlson", "AnnmarieConnor", "MoniqueWitte", "MikelWhitlock",
    "VincentAmaya", "KeiraQuiroz", "EllisCreech", "DionteCreel", "NicholeFoust",
    "ManuelConnell", "LourdesElam", "LincolnFrederick", "AlisaCash", "LucilleGriffiths"
)

$wsLogin = $wb.Worksheets.Item("login")
for ($idx = 0; $idx -lt $namesGroup1.Count; $idx++) {
    $row = $idx + 2
    $base = $namesGroup1[$idx]
    $newName = $base + "22"
    $newEmail = $base + "22@gmail.com"
    $wsLogin.Range("G$row").Value = $newName
    $wsLogin.Range("H$row").Value = $newName
    $wsLogin.Range("I$row").Value = $newEmail
}

# ---- "order": literal DonnellJernigan26 / ...26@gmail.com -> ...27 ----
$namesGroup2 = @(
    "DonnellJernigan", "MalikOtoole", "AlanCaudill", "AdanApplegate", "AiyanaWhitworth",
    "MercedezBrien", "DuaneHager", "LorenBell", "GeraldHiller", "DeionBranch",
    "DakotaHalstead", "ElliottFurman", "MiltonCamp", "DawnChester", "ZacheryPetrie",
    "EstebanAngel", "JimmyBlankenship", "AllysaGrice", "AugustineYoo", "BrandiSouthard"
)

$wsOrder = $wb.Worksheets.Item("order")
for ($idx = 0; $idx -lt $namesGroup2.Count; $idx++) {
    $row = $idx + 2
    $base = $namesGroup2[$idx]
    $newName = $base + "27"
    $newEmail = $base + "27@gmail.com"
    $wsOrder.Range("R$row").Value = $newName
    $wsOrder.Range("S$row").Value = $newName
    $wsOrder.Range("T$row").Value = $newEmail
}
